$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New task "AI Stuff" added in row 23 (with Estimated Hours 5, Day 1 burndown 3)
$ws.Range("A23").Value = "AI Stuff"
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = 3

# Gamepad Input (row 10) - Day 1 burndown updated from 0 to 1
$ws.Range("C10").Value = 1

# Update the active selection to reflect where the user ended up working
$ws.Range("C11").Select()
